$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Editorial Boards" -> "Guest editorial team"
$ws.Range("D2").Value = "Guest editorial team"

# Translate the "Editado por..." credit line to English, and swap "y" for "&"
$ws.Range("E3").Value = "Edited by Juan David Leongómez, Katarzyna Pisanski, David Reby, Disa Sauter, Nadine Lavan, Marcus Perlman & Jaroslava Varella Valentova"

# Reorder "Perfil \href{...}{Loop}" -> "\href{...}{Loop} profile"
$ws.Range("E5").Value = "\href{https://loop.frontiersin.org/people/438954/overview}{Loop} profile"

# Update the active selection to match the final state of the workbook
$ws.Range("E19").Select()
